$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SYNTHESIS")

# Update the ALU selector C2 (hex-ish code) 0 -> 4
$ws.Range("C2").Value = 4

# Select ALU operation "B+1" in H2 (drives G2 formula)
$ws.Range("H2").Value = "B+1"

# Select MEM operation "fetch" in R2 (was "write")
$ws.Range("R2").Value = "fetch"

# Select C bus destination "PC" in U2 (was "TOS")
$ws.Range("U2").Value = "PC"

# Set the JAM bit O3 to 1 (manual, non-formula cell)
$ws.Range("O3").Value = 1

# Add a new row 6 with a merged, empty separator cell B6:P6
$ws.Rows.Item(6).RowHeight = 18.75
$row6 = $ws.Range("B6:P6")
$row6.Merge()
$row6.Font.Name = "Consolas"
$row6.Font.Size = 14
$row6.HorizontalAlignment = -4108
$row6.NumberFormat = "0"
$row6.Borders.Item(8).LineStyle = 1
$row6.Borders.Item(8).Weight = 1
$row6.Borders.Item(8).ColorIndex = 8

# Update selection to match target sheet view state
$ws.Range("B4:P4").Select()
